$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows before row 371, shifting the existing rows (371..392) down to (373..394).
$ws.Rows.Item(371).Insert()
$ws.Rows.Item(371).Insert()

# New row 371
$ws.Cells.Item(371, 1).Value = 9
$ws.Cells.Item(371, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(371, 3).Value = "Metropolitana"
$ws.Cells.Item(371, 4).Value = 44516
$ws.Cells.Item(371, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(371, 5).Value = 13
$ws.Cells.Item(371, 6).Value = 100112040
$ws.Cells.Item(371, 7).Value = "Cilantro"
$ws.Cells.Item(371, 8).Value = "Sin especificar"
$ws.Cells.Item(371, 9).Value = "Primera"
$ws.Cells.Item(371, 10).Value = 43
$ws.Cells.Item(371, 11).Value = 5000
$ws.Cells.Item(371, 12).Value = 5000
$ws.Cells.Item(371, 13).Value = 5000
$ws.Cells.Item(371, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(371, 15).Value = "Región Metropolitana"
$ws.Cells.Item(371, 16).Value = 139
$ws.Cells.Item(371, 17).Value = 36
$ws.Cells.Item(371, 18).Value = "Hortaliza"

# New row 372
$ws.Cells.Item(372, 1).Value = 9
$ws.Cells.Item(372, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(372, 3).Value = "Metropolitana"
$ws.Cells.Item(372, 4).Value = 44516
$ws.Cells.Item(372, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(372, 5).Value = 13
$ws.Cells.Item(372, 6).Value = 100112040
$ws.Cells.Item(372, 7).Value = "Cilantro"
$ws.Cells.Item(372, 8).Value = "Sin especificar"
$ws.Cells.Item(372, 9).Value = "Primera"
$ws.Cells.Item(372, 10).Value = 160
$ws.Cells.Item(372, 11).Value = 12000
$ws.Cells.Item(372, 12).Value = 14000
$ws.Cells.Item(372, 13).Value = 13000
$ws.Cells.Item(372, 14).Value = "$/docena de atados"
$ws.Cells.Item(372, 15).Value = "Región Metropolitana"
$ws.Cells.Item(372, 16).Value = 4333
$ws.Cells.Item(372, 17).Value = 3
$ws.Cells.Item(372, 18).Value = "Hortaliza"
